$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab4")

# --- Data corrections (Lab 4 - Entrega Final) ---
# Insertion Sort timings for ARRAYLIST (rows 2-11) and LINKED_LIST (rows 15-24) tables
$ws.Range("B5").Value = 43601.2
$ws.Range("B6").Value = 179375.1
$ws.Range("B16").Value = 334934.56

# --- View state ---
$win = $excel.ActiveWindow
$win.Zoom = 125
[void]$ws.Range("F23").Select()
